$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 10 content
$ws.Range("A10").Value = "10. Fast & Slow Pointers"
$ws.Range("B10").Value = "Has Cycle"
$ws.Range("C10").Value = "Go throw the linked list with fast & slow pointer. If they both point to the same node then there is a cycle"

# Copy formatting from row 8 (an existing similarly-styled row) to row 10
$ws.Range("A8:C8").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122)

# Set row height for the new row
$ws.Rows.Item(10).RowHeight = 36.75

# Remove the fill/highlight color from cells that used the "theme4" colored style (C5,C6,C7,C9)
$ws.Range("C5").Interior.ColorIndex = -4142
$ws.Range("C6").Interior.ColorIndex = -4142
$ws.Range("C7").Interior.ColorIndex = -4142
$ws.Range("C9").Interior.ColorIndex = -4142

# Update the selection to A12
$ws.Range("A12").Select()
